# ---------------------------------------------------------------------------
# Default language + Resource update
#
#  1) Title paragraph ("TEMPLATE_TITLE"): explicit 36pt (sz/szCs=72 half-pts)
#     run size, set on both the paragraph mark run-properties and the run
#     itself (so Word no longer relies on the Title style's 45pt default).
#  2) "inetum.world" (footer-ish brand line): wrap it in spell-check
#     proofErr bookmarks, as Word does for words it does not recognise.
#  3) The legal/classification line: split it up so "Inetum", "Ref" and
#     "external" are each wrapped in their own proofErr spellStart/spellEnd
#     pair, matching what Word's spell checker would tag as unknown words.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1) Title: force explicit 36pt (sz/szCs = 72) on the paragraph mark and
#        on the run that holds "TEMPLATE_TITLE". ----------------------------

$titleRng = $d.Content
$null = $titleRng.Find.Execute("TEMPLATE_TITLE", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$titlePara = $titleRng.Paragraphs(1)
$titlePara.Range.Font.Size = 36
$titlePara.Range.Font.SizeBi = 36

# --- 2) "inetum.world": wrap the whole run in spellStart/spellEnd proofErr
#        marks, preserving its existing run formatting exactly. -------------

$worldRng = $d.Content
$null = $worldRng.Find.Execute("inetum.world", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$worldTarget = $d.Range($worldRng.Start, $worldRng.End)

$worldXml = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Arial"/><w:b/><w:color w:val="EE4641"/><w:sz w:val="16"/></w:rPr><w:t>inetum.world</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$worldTarget.InsertXML($worldXml)

# --- 3) Footer legal line: split "...Inetum | Ref. ..." and
#        "Classification: external" so the proper-noun / abbreviation /
#        English word each get their own proofErr spellStart/spellEnd. ------

$lineRng = $d.Content
$null = $lineRng.Find.Execute("12/08/2022", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$linePara = $lineRng.Paragraphs(1)
$lineTarget = $d.Range($linePara.Range.Start, $linePara.Range.End)

$lineXml = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r w:rsidRPr="0096307A"><w:rPr><w:rFonts w:ascii="Arial MT" w:hAnsi="Arial MT"/><w:sz w:val="14"/><w:lang w:val="fr-BE"/></w:rPr><w:t xml:space="preserve">12/08/2022 | &#169;2022 </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Arial MT" w:hAnsi="Arial MT"/><w:sz w:val="14"/><w:lang w:val="fr-BE"/></w:rPr><w:t>Inetum</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Arial MT" w:hAnsi="Arial MT"/><w:sz w:val="14"/><w:lang w:val="fr-BE"/></w:rPr><w:t xml:space="preserve"> | </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Arial MT" w:hAnsi="Arial MT"/><w:sz w:val="14"/><w:lang w:val="fr-BE"/></w:rPr><w:t>Ref</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Arial MT" w:hAnsi="Arial MT"/><w:sz w:val="14"/><w:lang w:val="fr-BE"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r>
<w:r w:rsidRPr="002D7FD2"><w:rPr><w:rFonts w:ascii="Arial MT" w:hAnsi="Arial MT"/><w:sz w:val="14"/><w:lang w:val="fr-BE"/></w:rPr><w:t>INTERNASSIGN v1.1</w:t></w:r>
<w:r w:rsidRPr="002D7FD2"><w:rPr><w:rFonts w:ascii="Arial MT" w:hAnsi="Arial MT"/><w:spacing w:val="-36"/><w:sz w:val="14"/><w:lang w:val="fr-BE"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r w:rsidRPr="002D7FD2"><w:rPr><w:rFonts w:ascii="Arial MT" w:hAnsi="Arial MT"/><w:sz w:val="14"/><w:lang w:val="fr-BE"/></w:rPr><w:t xml:space="preserve">Classification: </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Arial MT" w:hAnsi="Arial MT"/><w:sz w:val="14"/><w:lang w:val="fr-BE"/></w:rPr><w:t>external</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r w:rsidRPr="002D7FD2"><w:rPr><w:rFonts w:ascii="Arial MT" w:hAnsi="Arial MT"/><w:spacing w:val="-1"/><w:sz w:val="14"/><w:lang w:val="fr-BE"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r w:rsidRPr="002D7FD2"><w:rPr><w:rFonts w:ascii="Arial MT" w:hAnsi="Arial MT"/><w:sz w:val="14"/><w:lang w:val="fr-BE"/></w:rPr><w:t>document</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$lineTarget.InsertXML($lineXml)

Write-Host "edits applied"
